# Regenerate save_data column G ("K") values for rows 2-64.
# This corresponds to the commit: "regen save_data to use K instead of
# Strike#, regen std/mean, calc and write s_vals" — the recomputed s_vals
# (K column) are written back into the sheet, row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for G2:G64, in row order (63 values).
$newG = @(1,1,1,3,1,4,0,1,1,1,4,3,1,1,1,0,4,1,2,1,1,2,2,1,3,3,1,2,1,3,1,1,3,1,1,0,3,1,3,0,0,2,0,0,0,0,3,1,0,0,1,1,1,0,0,1,1,1,1,0,2,1,2)

$startRow = 2
for ($i = 0; $i -lt $newG.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newG[$i]
}
